# Logged Week 15 and simulated Week 16
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append newly logged per-play yardage samples to the four
# running lists (R/P, OFF/DEF).
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value2 + " 6 0 2 4 5 5 6 5 1 -1 5 8 4 3 5 22 2 12 1 8 2 4 1 1 3"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value2 + " 4 2 4 4 3 5 -1 5 2 2 1 3 2 6 7 2 2 5 6 14 4 5 7 1 0 -1 7 7 13 1 -2 10 1 2"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value2 + " 6 2 6 17 2 2 9 8 4 19 3 56 4 4 15"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value2 + " 5 8 7 28 -1 10 3 11 -1 8 9 11 8 25 2 24 24 5"

# ---------------------------------------------------------------------------
# OFF sheet: season offensive totals through the newly logged weeks.
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 190
$wsOFF.Range("F2").Value = 50
$wsOFF.Range("G2").Value = 51
$wsOFF.Range("H2").Value = 2
$wsOFF.Range("I2").Value = 4
$wsOFF.Range("J2").Value = 23
$wsOFF.Range("N2").Value = 16

$wsOFF.Range("C3").Value = 136
$wsOFF.Range("E3").Value = 26
$wsOFF.Range("F3").Value = 88
$wsOFF.Range("G3").Value = 32
$wsOFF.Range("H3").Value = 31
$wsOFF.Range("I3").Value = 46
$wsOFF.Range("J3").Value = 50
$wsOFF.Range("L3").Value = 201
$wsOFF.Range("M3").Value = 133
$wsOFF.Range("Q3").Value = 397

# ---------------------------------------------------------------------------
# DEF sheet: season defensive totals through the newly logged weeks.
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value = 154
$wsDEF.Range("F2").Value = 50
$wsDEF.Range("G2").Value = 39
$wsDEF.Range("J2").Value = 23
$wsDEF.Range("N2").Value = 20
$wsDEF.Range("O2").Value = 15

$wsDEF.Range("B3").Value = 9
$wsDEF.Range("C3").Value = 184
$wsDEF.Range("E3").Value = 37
$wsDEF.Range("F3").Value = 96
$wsDEF.Range("G3").Value = 34
$wsDEF.Range("H3").Value = 39
$wsDEF.Range("I3").Value = 54
$wsDEF.Range("J3").Value = 47
$wsDEF.Range("L3").Value = 276
$wsDEF.Range("M3").Value = 178
$wsDEF.Range("Q3").Value = 465

# ---------------------------------------------------------------------------
# ST sheet: special teams totals + newly logged per-event samples.
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 77
$wsST.Range("D2").Value = 56
$wsST.Range("H2").Value = 4
$wsST.Range("J2").Value = 49
$wsST.Range("K2").Value = 47
$wsST.Range("N2").Value = 20
$wsST.Range("O2").Value = 16

$wsST.Range("B3").Value = 51

$wsST.Range("D3").Value = $wsST.Range("D3").Value2 + " 57 57 41 61 45 56 51"
$wsST.Range("D4").Value = $wsST.Range("D4").Value2 + " 0 0 0 0 0 12 18"
$wsST.Range("D5").Value = $wsST.Range("D5").Value2 + " 12 0 0 5 0 5"
$wsST.Range("B6").Value = $wsST.Range("B6").Value2 + " 23"

# ---------------------------------------------------------------------------
# TURNS sheet: season turnover totals through the newly logged weeks.
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("D3").Value = 7
$wsTURNS.Range("E3").Value = 9

# ---------------------------------------------------------------------------
# PEN sheet: season penalty totals through the newly logged weeks.
# ---------------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value = 10
$wsPEN.Range("D2").Value = 5
$wsPEN.Range("B3").Value = 14
$wsPEN.Range("D4").Value = 6
